# Add new columns (area, houseno, locationLink) with a header row and one
# data row, matching the uploaded customers.xlsx revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "area"
$ws.Range("J1").Value = "houseno"
$ws.Range("K1").Value = "locationLink"

# --- Data row (row 2) ---
$ws.Range("I2").Value = "Al Rayyan"
$ws.Range("J2").Value = 54
$ws.Range("K2").Value = "https://www.google.com/maps/@25.2960247,51.5101904,14.75z"

# Turn the location link into a real hyperlink (applies the built-in
# "Hyperlink" style: blue, underlined text).
$ws.Hyperlinks.Add($ws.Range("K2"), "https://www.google.com/maps/@25.2960247,51.5101904,14.75z")

# Auto-size the new "area" and "locationLink" columns to fit their content,
# same as the original author's columns.
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(11).AutoFit()

# Leave the selection where the author left it before saving.
[void]$ws.Range("E11").Select()
